$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "41.524.79"
$ws.Range("E2").Value = "  +0.14%  "
$ws.Range("D3").Value = "2.470.75"
$ws.Range("E3").Value = "  -0.54%  "
$ws.Range("E4").Value = "  -0.21%  "
$ws.Range("D5").Value = "'314.74"
$ws.Range("E5").Value = "  +0.53%  "
$ws.Range("D6").Value = "'92.16"
$ws.Range("E6").Value = "  -2.68%  "
$ws.Range("D7").Value = "'0.548"
$ws.Range("E7").Value = "  -0.15%  "
$ws.Range("E8").Value = "  -0.21%  "
$ws.Range("D9").Value = "'0.516"
$ws.Range("E9").Value = "  +3.21%  "
$ws.Range("E10").Value = "  -4.00%  "
$ws.Range("E11").Value = "  +1.14%  "
$ws.Range("E12").Value = "  +0.70%  "
$ws.Range("E13").Value = "  -0.64%  "
$ws.Range("E14").Value = "  -2.02%  "
$ws.Range("D15").Value = "'15.97"
$ws.Range("E15").Value = "  +2.68%  "
$ws.Range("D16").Value = "2.459.38"
$ws.Range("E16").Value = "  -1.27%  "
$ws.Range("D17").Value = "'0.771"
$ws.Range("E17").Value = "  -2.93%  "
$ws.Range("D18").Value = "41.517.99"
$ws.Range("E18").Value = "  +0.24%  "
$ws.Range("E19").Value = "  +2.50%  "
$ws.Range("D20").Value = "0.0₃0947"
$ws.Range("E20").Value = "  +2.50%  "
$ws.Range("D21").Value = "'71.34"
$ws.Range("E21").Value = "  +3.42%  "
$ws.Range("D22").Value = "'11.10"
$ws.Range("E22").Value = "  -1.45%  "
$ws.Range("D23").Value = "'235.85"
$ws.Range("E23").Value = "  -0.64%  "
$ws.Range("E24").Value = "  -1.43%  "
$ws.Range("E25").Value = "  -0.09%  "
$ws.Range("D26").Value = "'1.90"
$ws.Range("E26").Value = "  -0.37%  "
$ws.Range("D27").Value = "'24.62"
$ws.Range("E27").Value = "  +1.82%  "
$ws.Range("E28").Value = "  -0.90%  "
$ws.Range("D29").Value = "'9.69"
$ws.Range("E29").Value = "  -0.65%  "
$ws.Range("D30").Value = "'35.38"
$ws.Range("E30").Value = "  -2.85%  "
$ws.Range("D31").Value = "'155.74"
$ws.Range("E32").Value = "  -0.87%  "
$ws.Range("E33").Value = "  +0.02%  "
$ws.Range("D34").Value = "'0.0759"
$ws.Range("E34").Value = "  +0.99%  "
$ws.Range("D35").Value = "'17.24"
$ws.Range("E35").Value = "  -4.81%  "
$ws.Range("D36").Value = "'2.88"
$ws.Range("E36").Value = "  -6.89%  "
$ws.Range("E37").Value = "  +1.23%  "
$ws.Range("E38").Value = "  -0.47%  "
$ws.Range("D39").Value = "'1.79"
$ws.Range("E39").Value = "  -4.61%  "
$ws.Range("D40").Value = "'2.22"
$ws.Range("E40").Value = "  -13.25%  "
$ws.Range("D41").Value = "'4.06"
$ws.Range("E41").Value = "  -3.16%  "
$ws.Range("E42").Value = "  -0.30%  "
$ws.Range("D43").Value = "1.944.79"
$ws.Range("E43").Value = "  -3.14%  "
$ws.Range("D45").Value = "'18.42"
$ws.Range("E45").Value = "  -5.78%  "
$ws.Range("E46").Value = "  -3.12%  "
$ws.Range("D47").Value = "'9.06"
$ws.Range("E47").Value = "  +2.87%  "
$ws.Range("D48").Value = "2.707.22"
$ws.Range("E48").Value = "  -0.94%  "
$ws.Range("D49").Value = "'96.99"
$ws.Range("E49").Value = "  -0.22%  "
$ws.Range("B50").Value = "MultiversX"
$ws.Range("C50").Value = "https://coinranking.com/coin/omwkOTglq+multiversx-egld"
$ws.Range("D50").Value = "'53.37"
$ws.Range("E50").Value = "  +4.50%  "
$ws.Range("B51").Value = "ordi"
$ws.Range("C51").Value = "https://coinranking.com/coin/j7-7vPrOi+ordi-ordi"
$ws.Range("D51").Value = "'67.14"
$ws.Range("E51").Value = "  -3.75%  "
